$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply corrected Monocytes_percent (column L) values for Buteo buteo rows 39-89
$ws.Range("L39").Value = 2
$ws.Range("L40").Value = 1
$ws.Range("L41").Value = 1
$ws.Range("L42").Value = 0
$ws.Range("L43").Value = 5
$ws.Range("L44").Value = 5
$ws.Range("L45").Value = 2
$ws.Range("L46").Value = 1
$ws.Range("L47").Value = 1
$ws.Range("L48").Value = 1
$ws.Range("L49").Value = 8
$ws.Range("L50").Value = 5
$ws.Range("L51").Value = 3
$ws.Range("L52").Value = 1
$ws.Range("L53").Value = 3
$ws.Range("L54").Value = 1
$ws.Range("L55").Value = 1
$ws.Range("L56").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("L58").Value = 1
$ws.Range("L59").Value = 0
$ws.Range("L60").Value = 1
$ws.Range("L61").Value = 3
$ws.Range("L62").Value = 3
$ws.Range("L63").Value = 4
$ws.Range("L64").Value = 3
$ws.Range("L65").Value = 4
$ws.Range("L66").Value = 6
$ws.Range("L67").Value = 2
$ws.Range("L68").Value = 4
$ws.Range("L69").Value = 6
$ws.Range("L70").Value = 3
$ws.Range("L71").Value = 4
$ws.Range("L72").Value = 4
$ws.Range("L73").Value = 4
$ws.Range("L74").Value = 2
$ws.Range("L75").Value = 4
$ws.Range("L76").Value = 1
$ws.Range("L77").Value = 1
$ws.Range("L78").Value = 2
$ws.Range("L79").Value = 3
$ws.Range("L80").Value = 2
$ws.Range("L81").Value = 2
$ws.Range("L82").Value = 7
$ws.Range("L83").Value = 4
$ws.Range("L84").Value = 6
$ws.Range("L85").Value = 7
$ws.Range("L86").Value = 1
$ws.Range("L87").Value = 3
$ws.Range("L88").Value = 3
$ws.Range("L89").Value = 4

# Apply the new font/style (fontId 1, cellXf 2) to the edited cells
$ws.Range("L39:L89").Font.Name = "Calibri"

# Update the view selection state to match the edit session
$ws.Range("L39:L89").Select()
